$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F14").Value = 12656
$ws1.Range("F16").Value = 5215
$ws1.Range("F17").Value = 5522

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 53

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F15").Value = 12656
$ws4.Range("F16").Value = 53
$ws4.Range("F19").Value = 5215
$ws4.Range("F20").Value = 5522
